# Update Government Revenue Accounting: carbon tax revenue now funds
# normal (regular) government spending instead of being split across
# deficit spending / household taxes / payroll taxes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Set Values Here")

# Row 8 = "carbon tax revenue" weights across B:F
# (Regular Spending, Deficit Spending, Household Taxes, Payroll Taxes, Corporate Taxes)
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0

# Reflect the edit's final UI state: the GRA-carbontax sheet's cursor
# moved from B5 to B3 (it is not the active tab)...
$gra = $wb.Worksheets.Item("GRA-carbontax")
[void]$gra.Select()
[void]$gra.Range("B3").Select()

# ...and the cursor on the About sheet moved off of the tab-selected
# cell as it was left (B31), while "Set Values Here" becomes the
# active/selected tab with C8 the active cell.
$about = $wb.Worksheets.Item("About")
[void]$about.Select()
[void]$about.Range("B31").Select()

[void]$ws.Select()
[void]$ws.Range("C8").Select()

[void]$wb.Application.Calculate()
